# Zadania.pptx edit: unhide the "Funkcja 1" solution variants (slides 3 & 4)
# and label all three variants (slides 3, 4, 5) as "Wersja A" / "Wersja B" / "Wersja C".

$p = $ppt.ActivePresentation

# --- Slide 3: "Wersja A" -------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.SlideShowTransition.Hidden = $False

$lbl3 = $s3.Shapes.Item(1).Duplicate()
$lbl3.Name = "pole tekstowe 9"
$lbl3.TextFrame.TextRange.Text = "Wersja A"
$lbl3.TextFrame.TextRange.Font.Bold = $True
$lbl3.TextFrame.WordWrap = $False
$lbl3.Left = 16.726456692913384
$lbl3.Top = 10.536614173228346
$lbl3.Width = 81.7051968503937
$lbl3.Height = 29.081259842519685

# --- Slide 4: "Wersja B" -------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.SlideShowTransition.Hidden = $False

$lbl4 = $s4.Shapes.Item(1).Duplicate()
$lbl4.Name = "pole tekstowe 9"
$lbl4.TextFrame.TextRange.Text = "Wersja B"
$lbl4.TextFrame.TextRange.Font.Bold = $True
$lbl4.TextFrame.WordWrap = $False
$lbl4.Left = 16.726456692913384
$lbl4.Top = 10.536614173228346
$lbl4.Width = 80.94787401574803
$lbl4.Height = 29.081259842519685

# --- Slide 5: "Wersja C" -------------------------------------------------
$s5 = $p.Slides.Item(5)

$lbl5 = $s5.Shapes.Item(1).Duplicate()
$lbl5.Name = "pole tekstowe 9"
$lbl5.TextFrame.TextRange.Text = "Wersja C"
$lbl5.TextFrame.TextRange.Font.Bold = $True
$lbl5.TextFrame.WordWrap = $False
$lbl5.Left = 16.726456692913384
$lbl5.Top = 10.536614173228346
$lbl5.Width = 80.31677165354331
$lbl5.Height = 29.081259842519685
